# ERP-435 - Split tribunal addresses onto multiple lines (Manchester & Glasgow)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the existing hyperlinks before we move cells around so stale
#    hyperlink ranges don't linger on the wrong cells.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Insert new rows so the single "address" row for each tribunal becomes
#    four rows (AddressLine1, AddressLine2, AddressLine3, Town) followed by
#    a dedicated PostCode row.
#    Manchester address used to live on row 3 (1 row) -> now rows 3-7 (5 rows)
#    so we need 4 new rows inserted right after row 3.
#    Glasgow address used to live on row 8 (post-insert) (1 row) -> now 4 rows
#    (AddressLine1, AddressLine2, Town, PostCode) so we need 3 new rows
#    inserted right after that row.
# ---------------------------------------------------------------------------

# Manchester: row 3 currently holds the combined address. Insert 4 rows below it.
$ws.Rows("4:7").Insert()

# Glasgow: after the Manchester insert, the combined Glasgow address row that
# used to be row 8 is now row 12. Insert 3 rows below it.
$ws.Rows("13:15").Insert()

# ---------------------------------------------------------------------------
# 3. Write the cell values (column A = field name, column B = value) for the
#    full, final layout.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "tribunalManchesterAddressLine1"
$ws.Range("B3").Value = "Manchester Employment Tribunal,"

$ws.Range("A4").Value = "tribunalManchesterAddressLine2"
$ws.Range("B4").Value = "Alexandra House,"

$ws.Range("A5").Value = "tribunalManchesterAddressLine3"
$ws.Range("B5").Value = "14-22 The Parsonage,"

$ws.Range("A6").Value = "tribunalManchesterTown"
$ws.Range("B6").Value = "Manchester,"

$ws.Range("A7").Value = "tribunalManchesterPostCode"
$ws.Range("B7").Value = "M3 2JA"

$ws.Range("A8").Value = "tribunalManchesterTelephone"
$ws.Range("B8").Value = "0161 833 6100"

$ws.Range("A9").Value = "tribunalManchesterFax"
$ws.Range("B9").Value = "0870 739 4433"

$ws.Range("A10").Value = "tribunalManchesterDX"
$ws.Range("B10").Value = "DX 743570"

$ws.Range("A11").Value = "tribunalManchesterEmail"
$ws.Range("B11").Value = "Manchesteret@justice.gov.uk"

$ws.Range("A12").Value = "tribunalGlasgowAddressLine1"
$ws.Range("B12").Value = "Eagle Building,"

$ws.Range("A13").Value = "tribunalGlasgowAddressLine2"
$ws.Range("B13").Value = "215 Bothwell Street,"

$ws.Range("A14").Value = "tribunalGlasgowTown"
$ws.Range("B14").Value = "Glasgow,"

$ws.Range("A15").Value = "tribunalGlasgowPostCode"
$ws.Range("B15").Value = "G2 7TS"

$ws.Range("A16").Value = "tribunalGlasgowTelephone"
$ws.Range("B16").Value = "0141 204 0730"

$ws.Range("A17").Value = "tribunalGlasgowFax"
$ws.Range("B17").Value = "01264 785 177"

$ws.Range("A18").Value = "tribunalGlasgowDX"
$ws.Range("B18").Value = "DX 580003"

$ws.Range("A19").Value = "tribunalGlasgowEmail"
$ws.Range("B19").Value = "glasgowet@justice.gov.uk"

# ---------------------------------------------------------------------------
# 4. Formatting.
#    Column A labels for the Manchester address block (rows 2-6) wrap,
#    matching the rest of the "FIELDS" column for that tribunal's block.
# ---------------------------------------------------------------------------
$ws.Range("A2:A6").WrapText = $true
$ws.Range("A7:A19").WrapText = $false

# Column B: only the address-line cells wrap; the town/postcode/contact
# values do not.
$ws.Range("B3:B6").WrapText = $false
$ws.Range("B7").WrapText = $true
$ws.Range("B8:B11").WrapText = $false
$ws.Range("B12:B13").WrapText = $true
$ws.Range("B14").WrapText = $false
$ws.Range("B15").WrapText = $true
$ws.Range("B16:B19").WrapText = $false

# Give the Glasgow "Town" value its own (plain Calibri 11) font, matching the
# distinct font used for that single cell in the source workbook.
$ws.Range("B14").Font.Name = "Calibri"
$ws.Range("B14").Font.Size = 11
$ws.Range("B14").Font.Bold = $false
$ws.Range("B14").Font.Italic = $false
$ws.Range("B14").Font.Color = 0

# The "H" helper column keeps its right-aligned placeholder style down to
# the last Manchester row (rows 1-7).
$ws.Range("H4:H7").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 5. Re-create the hyperlinks on the e-mail cells, which have now moved down
#    to rows 11 (Manchester) and 19 (Glasgow). Adding a hyperlink auto-applies
#    the "Hyperlink" (underlined/blue) font, so restore the plain font the
#    source workbook actually used for these cells.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:Manchesteret@justice.gov.uk", "", "", "Manchesteret@justice.gov.uk")
$ws.Range("B11").Font.Color = 0
$ws.Range("B11").Font.Underline = 0
$ws.Range("B11").Font.Name = "Calibri"
$ws.Range("B11").Font.Size = 11

$ws.Hyperlinks.Add($ws.Range("B19"), "mailto:glasgowet@justice.gov.uk", "", "", "glasgowet@justice.gov.uk")
$ws.Range("B19").Font.Color = 0
$ws.Range("B19").Font.Underline = 0
$ws.Range("B19").Font.Name = "Calibri"
$ws.Range("B19").Font.Size = 11

# ---------------------------------------------------------------------------
# 6. Move the sheet selection, mirroring the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("A6").Select()
